# Daily update at 8 AM UTC
# Appends the next day's row to the bottom of the "Wins Over Time" tracking
# sheet, moving the special "last row" date style down onto the newly
# appended row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the current last data row (the row that currently carries the
# distinctive "last row" date-only number format).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

# That row is no longer the final entry once we append a new one, so give it
# the regular date/time number format used by every other data row.
$ws.Cells.Item($lastRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the next day's data: date = previous date + 1, plus the day's counts.
$prevDate = $ws.Cells.Item($lastRow, 1).Value2
$ws.Cells.Item($newRow, 1).Value2 = $prevDate + 1
$ws.Cells.Item($newRow, 2).Value = 104
$ws.Cells.Item($newRow, 3).Value = 106
$ws.Cells.Item($newRow, 4).Value = 103

# The newly appended row becomes the new "last row", so it takes on the
# date-only number format that used to mark the final entry.
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD"
